# Update the "Parallel_Use" (G2) value in the Setup sheet from 2 to 8,
# matching the last experiment_manager parameters, and leave the
# selection on the edited cell (as it would be after typing into it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("G2")
$cell.Select()
$cell.Value = 8
